$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 21666.666
$ws.Range("L21").Value = 21666.666
$ws.Range("N21").Value = -22602.666
# Row 23
$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 21666.666
$ws.Range("L23").Value = 21666.666
$ws.Range("N23").Value = -22134.666
# Row 29
$ws.Range("H29").Value = 725
$ws.Range("J29").Value = 750
$ws.Range("L29").Value = 2250
$ws.Range("N29").Value = -2812
# Row 43
$ws.Range("H43").Value = 7962.5
$ws.Range("I43").Value = 7925
$ws.Range("J43").Value = 8000
$ws.Range("K43").Value = 7925
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = -7856
$ws.Range("N43").Value = -8138
# Row 64
$ws.Range("H64").Value = 6916.6665
$ws.Range("J64").Value = 9500
$ws.Range("L64").Value = 9500
$ws.Range("N64").Value = -9996
# Row 67
$ws.Range("H67").Value = 6916.6665
$ws.Range("J67").Value = 9500
$ws.Range("L67").Value = 9500
$ws.Range("N67").Value = -11216
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2428.4285
$ws.Range("I45").Value = 2599.8
$ws.Range("K45").Value = 2599.8
$ws.Range("M45").Value = -2222.8
# Row 63
$ws.Range("H63").Value = 12000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66
$ws.Range("H66").Value = 12000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 92
$ws.Range("H92").Value = 160183.33
$ws.Range("J92").Value = 160183.33
$ws.Range("L92").Value = 160183.33
$ws.Range("N92").Value = -165175.33
# Row 96
$ws.Range("H96").Value = 28498.75
$ws.Range("J96").Value = 28498.75
$ws.Range("L96").Value = 28498.75
$ws.Range("N96").Value = -33990.75
# Row 106
$ws.Range("H106").Value = 7500
$ws.Range("J106").Value = 7500
$ws.Range("L106").Value = 7500
$ws.Range("N106").Value = -10024
# Row 110
$ws.Range("H110").Value = 83334060
$ws.Range("I110").Value = 1099.5
$ws.Range("J110").Value = 250000000
$ws.Range("K110").Value = 1099.5
$ws.Range("L110").Value = 250000000
$ws.Range("M110").Value = 945.5
$ws.Range("N110").Value = -250004090
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 132
$ws.Range("H132").Value = 2491.6667
$ws.Range("I132").Value = 2485
$ws.Range("J132").Value = 2495
$ws.Range("K132").Value = 7455
$ws.Range("L132").Value = 7485
$ws.Range("N132").Value = -12545
$ws.Range("M132").Value = -4925
# Row 139
$ws.Range("H139").Value = 99473.75
$ws.Range("J139").Value = 99473.75
$ws.Range("L139").Value = 99473.75
$ws.Range("N139").Value = -109753.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 909
$ws.Range("I64").Value = 909
$ws.Range("K64").Value = 909
$ws.Range("M64").Value = -684
# Row 67
$ws.Range("H67").Value = 909
$ws.Range("I67").Value = 909
$ws.Range("K67").Value = 909
$ws.Range("M67").Value = -129
# Row 99
$ws.Range("H99").Value = 3383.625
$ws.Range("I99").Value = 2084.6924
$ws.Range("K99").Value = 2084.6924
$ws.Range("M99").Value = -586.6923999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 829.8889
$ws.Range("I22").Value = 846.125
$ws.Range("K22").Value = 846.125
$ws.Range("M22").Value = -496.125
# Row 62
$ws.Range("H62").Value = 3250
$ws.Range("I62").Value = 3250
$ws.Range("K62").Value = 3250
$ws.Range("M62").Value = -2626
# Row 65
$ws.Range("H65").Value = 3250
$ws.Range("I65").Value = 3250
$ws.Range("K65").Value = 16250
$ws.Range("M65").Value = -13130
# Row 96
$ws.Range("H96").Value = 36274.668
$ws.Range("J96").Value = 36274.668
$ws.Range("L96").Value = 36274.668
$ws.Range("N96").Value = -41766.668

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1672.8572
$ws.Range("I4").Value = 499.23077
$ws.Range("J4").Value = 2999.5652
$ws.Range("K4").Value = 1497.69231
$ws.Range("L4").Value = 8998.695599999999
$ws.Range("M4").Value = -1385.69231
$ws.Range("N4").Value = -9222.695599999999
# Row 32
$ws.Range("H32").Value = 5787
$ws.Range("I32").Value = 795
$ws.Range("K32").Value = 2385
$ws.Range("M32").Value = -2102
# Row 136
$ws.Range("H136").Value = 1029.6666
$ws.Range("I136").Value = 1029.6666
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3088.9998
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 2011.0002
$ws.Range("N136").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 150.16667
$ws.Range("I2").Value = 156.29411
$ws.Range("K2").Value = 156.29411
$ws.Range("M2").Value = -43.29410999999999
# Row 63
$ws.Range("H63").Value = 55000
$ws.Range("J63").Value = 55000
$ws.Range("L63").Value = 55000
$ws.Range("N63").Value = -56372
# Row 66
$ws.Range("H66").Value = 55000
$ws.Range("J66").Value = 55000
$ws.Range("L66").Value = 165000
$ws.Range("N66").Value = -171864
# Row 92
$ws.Range("H92").Value = 12406.667
$ws.Range("J92").Value = 14688
$ws.Range("L92").Value = 14688
$ws.Range("N92").Value = -18432
# Row 132
$ws.Range("H132").Value = 5012
$ws.Range("I132").Value = 5012
$ws.Range("K132").Value = 15036
$ws.Range("M132").Value = -12506

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 944
$ws.Range("I16").Value = 944
$ws.Range("K16").Value = 944
$ws.Range("M16").Value = -774
# Row 68
$ws.Range("H68").Value = 6055.4443
$ws.Range("I68").Value = 4083.1667
$ws.Range("K68").Value = 4083.1667
$ws.Range("M68").Value = -3334.1667
# Row 71
$ws.Range("H71").Value = 6055.4443
$ws.Range("I71").Value = 4083.1667
$ws.Range("K71").Value = 20415.8335
$ws.Range("M71").Value = -16671.8335
